$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C", "D", "E", "F", "G", "I", "K", "L", "M", "N", "O")

$data = @{}
$data[2] = @(9.778152741487753, 5.546148107585794, 12.4848175053883, 27.44739673838564, 3.651785335344727, 25.32755383773117, 13.68965849097139, 9.541650222979287, 16.55137498137744, 18.04841408333072, 24.5579473949477)
$data[3] = @(9.752335347992043, 5.51391437142799, 12.50169040778178, 27.46674358416147, 3.653862695487962, 25.38136430671656, 13.30117255470987, 9.56570237314407, 16.40082239024794, 18.0931657096762, 24.62353342766957)
$data[4] = @(9.73821545968217, 5.493761155960094, 12.51416138523469, 27.48629850380663, 3.655206371427738, 25.42027808278433, 13.05831602999087, 9.581588236039359, 16.31013318831621, 18.12243248604794, 24.66960438277025)
$data[5] = @(9.732901315474333, 5.48546024391575, 12.51977418285101, 27.49619436681756, 3.655771125900995, 25.43760952160191, 12.95840796502116, 9.588343197559443, 16.27365078924971, 18.13480964145922, 24.68983338371248)
$data[6] = @(9.73204558739447, 5.484076626146912, 12.5207382377201, 27.49795384595274, 3.655865943116327, 25.44057630217275, 12.94176584169644, 9.589481854404793, 16.26762253036891, 18.13689210474011, 24.69328013406082)
$data[7] = @(9.738142005079402, 5.493649561096814, 12.51423493243694, 27.48642416515251, 3.65521391820788, 25.4205058580633, 13.05697225023959, 9.58167819629511, 16.30963920972417, 18.1225975827043, 24.66987131352826)
$data[8] = @(9.768894088636925, 5.535110518016541, 12.49019715893239, 27.45247287273303, 3.652487493394716, 25.34488667865988, 13.55669038812072, 9.549711583213954, 16.49912170353333, 18.06347360924941, 24.579355183058)
$data[9] = @(9.842737764286005, 5.613433846095511, 12.45980740360474, 27.44688968148888, 3.647679399429286, 25.24333811872287, 14.49642725490341, 9.495881578645136, 16.88306176462532, 17.96169075370775, 24.44804770576543)
$data[10] = @(9.904946891097312, 5.66901251912207, 12.44768404004192, 27.48002559456026, 3.644471670714005, 25.1973864426861, 15.15515656420622, 9.461714838509657, 17.17060274772883, 17.8954927938444, 24.37995063680686)
$data[11] = @(9.934904053710216, 5.693838860727237, 12.44438088622389, 27.50317215562929, 3.643082179454187, 25.18273060914355, 15.4466323293196, 9.447336645712141, 17.30215085023303, 17.86723086477107, 24.3551722493715)
$data[12] = @(9.946479897253841, 5.703171636809082, 12.4434475813614, 27.51309470377905, 3.642565985311223, 25.17808064812712, 15.55573368573024, 9.44205918086721, 17.35203414907289, 17.8567943579331, 24.34668344564882)
$data[13] = @(9.943976636533044, 5.701164743025307, 12.44363447295787, 27.51090629583309, 3.642676714055357, 25.17904205955175, 15.53229487594444, 9.443188342235672, 17.34128842849111, 17.85903024025255, 24.34847184969145)
$data[14] = @(9.935851800852014, 5.69460806074341, 12.44429774244013, 27.50396531682219, 3.643039512200168, 25.18233001232405, 15.45563422418775, 9.446899114800287, 17.30625361607377, 17.86636692609419, 24.35445593395991)
$data[15] = @(9.930905076619224, 5.690582911959424, 12.4447453476887, 27.49986437717068, 3.643263034245658, 25.18446120524253, 15.40850869689602, 9.449193844481195, 17.28480162592015, 17.87089543828012, 24.35823789037211)
$data[16] = @(9.90302191715706, 5.667380634015444, 12.44794437415928, 27.47867509603624, 3.644563875934569, 25.1984700929925, 15.13593471727021, 9.46267790147016, 17.1620174805453, 17.89737698777879, 24.38169495455057)
$data[17] = @(9.88633647881926, 5.65302785635344, 12.450473038948, 27.46774178772973, 3.645379721812243, 25.20866549024676, 14.96655103440636, 9.471248029729566, 17.08685680045832, 17.91409641595508, 24.39767498754642)
$data[18] = @(9.876895928955708, 5.644729696732386, 12.45213568741161, 27.46221338069851, 3.645855540163116, 25.2151176287313, 14.86836094657443, 9.476286954104879, 17.04369869084426, 17.92388733003738, 24.40744964427281)
$data[19] = @(9.873726593848991, 5.64191281211295, 12.45273440887901, 27.46047219552537, 3.646017773212238, 25.21740315227718, 14.8349871943908, 9.478011879595662, 17.02909968450288, 17.92723232789001, 24.41085927794347)
$data[20] = @(9.888096525288988, 5.654560185915466, 12.45018231094241, 27.46882700581806, 3.645292194494438, 25.20751930334711, 14.98466213498303, 9.470324381828558, 17.09485057480331, 17.91229856518528, 24.39591348842325)
$data[21] = @(9.93823202824381, 5.696535795968496, 12.44409431137302, 27.50597267144082, 3.642932679238087, 25.18133982896572, 15.47818662591682, 9.445804633547288, 17.31654262266877, 17.86420475820393, 24.35267397373357)
$data[22] = @(9.972345771451943, 5.723569220231779, 12.44196590839863, 27.5369935090854, 3.641448724835154, 25.16947559304815, 15.79326792689518, 9.430754300102496, 17.46181526719505, 17.83432094263928, 24.32962776291969)
$data[23] = @(9.95401763950623, 5.709178480172469, 12.44293276908875, 27.51982149012652, 3.642235437085036, 25.17532742900649, 15.62581592925679, 9.438697822686875, 17.38425789722972, 17.85012903110264, 24.34145011627554)
$data[24] = @(9.887300333917405, 5.65386756451648, 12.45031309838252, 27.46833401952467, 3.645331744478695, 25.20803565464506, 14.97647661922482, 9.470741614644698, 17.0912364221358, 17.91311081684302, 24.39670803226644)
$data[25] = @(9.82134087180274, 5.592579812849409, 12.46623578982444, 27.44185978365287, 3.648922833671961, 25.26578748671079, 14.24728956378831, 9.50949764195034, 16.77808515138567, 17.98771511735332, 24.47860120811982)

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = $cols[$i] + $r
        $ws.Range($addr).Value = $rowVals[$i]
    }
}